$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates reflecting refreshed crypto price/volume snapshot.
# Values that parse as plain numbers are written with a leading apostrophe
# so Excel keeps them as literal text (matching the original inline-string cells)
# instead of silently converting them to numbers (e.g. "6.50" -> 6.5).

$ws.Range('D2').Value = '59.060.40'
$ws.Range('E2').Value = '  +0.06%  '

$ws.Range('D3').Value = '2.516.45'
$ws.Range('E3').Value = '  -0.08%  '

$ws.Range('E4').Value = '  +0.08%  '

$ws.Range('D5').Value = '''534.52'
$ws.Range('E5').Value = '  +0.01%  '

$ws.Range('D6').Value = '''136.22'
$ws.Range('E6').Value = '  -0.34%  '

$ws.Range('E7').Value = '  -0.01%  '

$ws.Range('E8').Value = '  +0.27%  '

$ws.Range('D9').Value = '''0.102'
$ws.Range('E9').Value = '  +0.78%  '

$ws.Range('D10').Value = '''0.157'
$ws.Range('E10').Value = '  -1.11%  '

$ws.Range('D11').Value = '''5.41'
$ws.Range('E11').Value = '  +1.46%  '

$ws.Range('D12').Value = '''0.347'
$ws.Range('E12').Value = '  +0.04%  '

$ws.Range('D13').Value = '2.967.06'
$ws.Range('E13').Value = '  +0.51%  '

$ws.Range('D14').Value = '58.979.06'
$ws.Range('E14').Value = '  +0.06%  '

$ws.Range('D15').Value = '''22.87'
$ws.Range('E15').Value = '  -1.34%  '

$ws.Range('E16').Value = '  -0.63%  '

$ws.Range('D17').Value = '2.508.77'
$ws.Range('E17').Value = '  -0.29%  '

$ws.Range('D18').Value = '''11.09'
$ws.Range('E18').Value = '  +0.63%  '

$ws.Range('D19').Value = '''4.26'
$ws.Range('E19').Value = '  +0.43%  '

$ws.Range('D20').Value = '''324.61'
$ws.Range('E20').Value = '  -0.01%  '

$ws.Range('E21').Value = '  -0.07%  '

$ws.Range('D22').Value = '''5.96'
$ws.Range('E22').Value = '  +1.71%  '

$ws.Range('D23').Value = '''65.28'
$ws.Range('E23').Value = '  +0.46%  '

$ws.Range('D24').Value = '''0.422'
$ws.Range('E24').Value = '  +0.61%  '

$ws.Range('E25').Value = '  -0.40%  '

$ws.Range('E26').Value = '  +0.11%  '

$ws.Range('D27').Value = '''7.57'
$ws.Range('E27').Value = '  -0.11%  '

$ws.Range('D28').Value = '0.0₃0768'
$ws.Range('E28').Value = '  -0.36%  '

$ws.Range('D29').Value = '''6.50'
$ws.Range('E29').Value = '  -3.81%  '

$ws.Range('E30').Value = '  -0.58%  '

$ws.Range('D31').Value = '''168.99'
$ws.Range('E31').Value = '  +0.61%  '

$ws.Range('E32').Value = '  +0.04%  '

$ws.Range('E33').Value = '  -3.05%  '

$ws.Range('E34').Value = '  -2.16%  '

$ws.Range('D35').Value = '''18.41'
$ws.Range('E35').Value = '  -0.63%  '

$ws.Range('D36').Value = '''4.07'
$ws.Range('E36').Value = '  -0.82%  '

$ws.Range('D37').Value = '''1.53'
$ws.Range('E37').Value = '  -2.41%  '

$ws.Range('B38').Value = 'SuiNetwork'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D38').Value = '''0.806'
$ws.Range('E38').Value = '  -2.38%  '

$ws.Range('B39').Value = 'Filecoin'
$ws.Range('C39').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D39').Value = '''3.59'
$ws.Range('E39').Value = '  -0.80%  '

$ws.Range('D40').Value = '''283.13'
$ws.Range('E40').Value = '  +1.30%  '

$ws.Range('B41').Value = 'FirstDigitalUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D41').Value = '''0.998'
$ws.Range('E41').Value = '  -0.14%  '

$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').Value = '''5.04'
$ws.Range('E42').Value = '  -3.10%  '

$ws.Range('D43').Value = '''0.606'
$ws.Range('E43').Value = '  +0.48%  '

$ws.Range('D44').Value = '''130.39'
$ws.Range('E44').Value = '  +1.77%  '

$ws.Range('D45').Value = '''10.93'
$ws.Range('E45').Value = '  +0.36%  '

$ws.Range('E46').Value = '  -0.43%  '

$ws.Range('D47').Value = '''0.0502'
$ws.Range('E47').Value = '  -2.31%  '

$ws.Range('E48').Value = '  -1.36%  '

$ws.Range('D49').Value = '''17.36'
$ws.Range('E49').Value = '  +0.29%  '

$ws.Range('D50').Value = '1.760.27'
$ws.Range('E50').Value = '  -0.68%  '

$ws.Range('D51').Value = '''0.983'
$ws.Range('E51').Value = '  -0.54%  '
